$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")

# Update values in rows 5-7 (labels stay the same, values change)
$ws1.Range("B5").Value = 62100
$ws1.Range("B6").Value = 69300
$ws1.Range("B7").Value = "01.01.2024"

# Update selection on Tabelle1
$ws1.Range("A2:B7").Select()

# Remove data validation from B2
$ws1.Range("B2").Validation.Delete()
